# "Modulos de actividades.xlsx" update
#
# Source data changes (per the commit diff):
#  - "base de dados" sheet: Normalizacao progress goes from "EP" (em progresso)
#    to "F" (feito); BD Fisica progress goes from blank to "ep" (em progresso).
#    The stray note in D5 ("folege") is cleared along with it.
#  - "front-end" sheet: the three "instituicoes da saude" page rows (24-26)
#    get their Interface status set to "ep" and their developer set to "Neima".
#  - The sheet also ends up protected, and the remembered selections for the
#    front-end and base-de-dados tabs change.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("front-end")
$ws3 = $wb.Worksheets.Item("base de dados")

# --- "base de dados" sheet -------------------------------------------------
$ws3.Range("C5").Value = "F"
$ws3.Range("D5").Value = ""
$ws3.Range("C6").Value = "ep"

# --- "front-end" sheet ------------------------------------------------------
$ws1.Range("E24").Value = "ep"
$ws1.Range("G24").Value = "Neima"
$ws1.Range("E25").Value = "ep"
$ws1.Range("G25").Value = "Neima"
$ws1.Range("E26").Value = "ep"
$ws1.Range("G26").Value = "Neima"

# --- remembered selections --------------------------------------------------
# "base de dados" keeps its own last-used selection without becoming the
# active tab.
$ws3.Activate()
$ws3.Range("E9").Select()

# "front-end" stays the active/selected tab, scrolled down to its last edit.
$ws1.Activate()
$ws1.Range("I26").Select()

# --- protect the front-end sheet -------------------------------------------
$ws1.Protect()
